# Added first test case; updated requirements table
#
# The "Requirements" sheet contains a Test-Case table (columns A:H).
# Column E is "TEST CASE". A first test case ("Test case 1") is being
# recorded against the requirement rows that call for test coverage.
# This replaces several now-obsolete one-off TEST CASE notes (which
# previously lived only as ad-hoc free text) with the single canonical
# placeholder string "Test case 1" used consistently across the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$testCaseRows = @(4, 6, 18, 19, 20, 21, 22, 24, 26, 27, 28, 29, 30, 31, 32, 36, 37, 38, 39, 40, 41)

foreach ($row in $testCaseRows) {
    $ws.Cells.Item($row, 5).Value = "Test case 1"
}
